# Bugfixed evaluation and simulated rt_data for components.
# Replaces the "YYYYQ4" text labels in column A with real date values
# (end-of-year / Q4 dates) formatted as "YYYY-MM-DD HH:MM:SS".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Serial date numbers (OLE Automation date) for 2004-12-31 .. 2024-12-31,
# one per row from A2 to A22.
$dates = @(38352, 38717, 39082, 39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735, 43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657)

for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $dates[$i]
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
